$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Right" count for Marking row (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Update "Right" count for Total row (B12): 57 -> 95
$ws.Range("B12").Value = 95

# Update displayed "Corr/total" marks text for Total row (E12): 56/84 -> 95/140
$ws.Range("E12").Value = "95/140"
